$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.804.51"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").Value = "1.557.55"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'205.22"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "'0.479"
$ws.Range("E6").Value = "  -1.83%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'21.50"
$ws.Range("E8").Value = "  -3.57%  "
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").Value = "'0.0584"
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").Value = "1.779.48"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "1.555.21"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "'3.71"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").Value = "'0.512"
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("D16").Value = "26.827.45"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").Value = "'61.27"
$ws.Range("E17").Value = "  -2.49%  "
$ws.Range("D18").Value = "'214.23"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'7.29"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "0.0₃0682"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'4.11"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "'9.12"
$ws.Range("E23").Value = "  -2.86%  "
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").Value = "'152.77"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").Value = "'6.56"
$ws.Range("E26").Value = "  -1.26%  "
$ws.Range("D27").Value = "'14.93"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("D32").Value = "'3.18"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").Value = "1.370.63"
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -2.99%  "
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").Value = "'0.924"
$ws.Range("E37").Value = "  -2.14%  "
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("D39").Value = "'0.524"
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("D40").Value = "'0.807"
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "'5.58"
$ws.Range("E42").Value = "  +4.21%  "
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("E45").Value = "  -1.46%  "
$ws.Range("D46").Value = "'63.33"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("E47").Value = "  -2.40%  "
$ws.Range("D48").Value = "1.693.24"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D49").Value = "'86.43"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("D50").Value = "'0.0509"
$ws.Range("E50").Value = "  +3.21%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0951"
$ws.Range("E51").Value = "  +0.18%  "
